# Update latest output (run 66)
$wb = $excel.ActiveWorkbook

# ---- Sheet "Schedule" (sheet1) ----
$wsSchedule = $wb.Worksheets.Item("Schedule")

$wsSchedule.Range("A4").Value = 46039.54166666666
$wsSchedule.Range("B4").Value = 46039.83333333334
$wsSchedule.Range("E4").Value = 260.3643315
$wsSchedule.Range("F4").Value = 9.839921825396827

$wsSchedule.Range("A5").Value = 46040.3125
$wsSchedule.Range("B5").Value = 46040.8125
$wsSchedule.Range("E5").Value = -29.91509624999999
$wsSchedule.Range("F5").Value = -0.6595038855820103

# ---- Sheet "Detailed" (sheet2) ----
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("E27").Value = "OFF"

$wsDetailed.Range("B40").Value = 0.01082

$wsDetailed.Range("B41").Value = 8.57385
$wsDetailed.Range("E41").Value = "ON"

$wsDetailed.Range("B42").Value = 57.3

$wsDetailed.Range("B43").Value = 57.3
$wsDetailed.Range("C43").Value = "historical"

$wsDetailed.Range("C44").Value = "historical"

$wsDetailed.Range("B45").Value = 56.98

$wsDetailed.Range("B47").Value = 57.06004

$wsDetailed.Range("B49").Value = 36.2

$wsDetailed.Range("B50").Value = 36.2

$wsDetailed.Range("B52").Value = 36.0603

$wsDetailed.Range("B53").Value = 52.28671

$wsDetailed.Range("B54").Value = 36.2

$wsDetailed.Range("B55").Value = 36.2

$wsDetailed.Range("B56").Value = 56.97996

$wsDetailed.Range("B59").Value = 56.98

$wsDetailed.Range("B60").Value = 57.06003

$wsDetailed.Range("B64").Value = 36.0595
$wsDetailed.Range("E64").Value = "OFF"

$wsDetailed.Range("B65").Value = 27.51631

$wsDetailed.Range("B67").Value = 21.8653

$wsDetailed.Range("B68").Value = -0.9349499999999999

$wsDetailed.Range("B69").Value = -5.01

$wsDetailed.Range("B70").Value = -3.6481

$wsDetailed.Range("B71").Value = 0.51

$wsDetailed.Range("B72").Value = 0.7

$wsDetailed.Range("B73").Value = 0.7

$wsDetailed.Range("B74").Value = -5.51

$wsDetailed.Range("B75").Value = -0.89546

$wsDetailed.Range("B76").Value = -5.51

$wsDetailed.Range("B77").Value = -5.59641

$wsDetailed.Range("B78").Value = -9.99

$wsDetailed.Range("B79").Value = -10.81131

$wsDetailed.Range("B80").Value = -13.30111

$wsDetailed.Range("B81").Value = -11.01

$wsDetailed.Range("B82").Value = -7.35773

$wsDetailed.Range("B83").Value = -7.29212

$wsDetailed.Range("B84").Value = -6.37938

$wsDetailed.Range("B85").Value = -0.62352

$wsDetailed.Range("B86").Value = -1.38933

$wsDetailed.Range("B87").Value = 0.00023

$wsDetailed.Range("B88").Value = 12.77543
$wsDetailed.Range("E88").Value = "ON"

$wsDetailed.Range("B89").Value = 44.0223

$wsDetailed.Range("B90").Value = 45.92752

$wsDetailed.Range("B91").Value = 45.17477

$wsDetailed.Range("B92").Value = 46.84454

$wsDetailed.Range("B94").Value = 56.03123
